# Weekly Fruta/Hortaliza update: insert a new daily-price record as the
# new row 98 (pushing the existing rows 98-179 down to 99-180), matching
# the latest "Segunda" grade Caramelo pineapple entry for the Vega
# Monumental Concepción market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 98; this shifts rows 98:179 down to 99:180
# and grows the sheet dimension from A1:T179 to A1:T180 automatically.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new weekly record.
$ws.Cells.Item(98, 1).Value  = 11
$ws.Cells.Item(98, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(98, 3).Value  = "Bíobío"
$ws.Cells.Item(98, 4).Value  = 44741
$ws.Cells.Item(98, 5).Value  = 8
$ws.Cells.Item(98, 6).Value  = "Fruta"
$ws.Cells.Item(98, 7).Value  = 100108
$ws.Cells.Item(98, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(98, 9).Value  = 100108005
$ws.Cells.Item(98, 10).Value = "Piña"
$ws.Cells.Item(98, 11).Value = "Caramelo"
$ws.Cells.Item(98, 12).Value = "Segunda"
$ws.Cells.Item(98, 13).Value = 200
$ws.Cells.Item(98, 14).Value = 17000
$ws.Cells.Item(98, 15).Value = 18000
$ws.Cells.Item(98, 16).Value = 17500
$ws.Cells.Item(98, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(98, 18).Value = "Ecuador"
$ws.Cells.Item(98, 19).Value = 1250
$ws.Cells.Item(98, 20).Value = 14
